# remove general journals from landing page gauges
#
# The gauges on the landing page were summing a "general journals" total
# that shouldn't be in the count, which was masking a missing timesheet
# entry. This appends the missing entry (row 64) to the Sheet1 log,
# continuing the running-total formula in column C so the totals in D2/F2
# (which the gauges read from) pick up the correct hours, and moves the
# active cell/viewport down to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New timesheet entry: date, hours, running total (previous total + hours)
$ws.Range("A64").Value = 45432
$ws.Range("B64").Value = 5.5
$ws.Range("C64").Formula = "=C63+B64"

# Match row 63's date formatting on the new date cell
[void]$ws.Range("A63").Copy()
$ws.Range("A64").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active selection / viewport to follow the newly added row
[void]$ws.Range("C64").Select()
$excel.ActiveWindow.ScrollRow = 36
